# Slide 3 ("Arquitectura por capas") has five numbered circles laid out
# top-to-bottom that originally read 1, 2, 3, 4, 5. The edit reverses the
# numbering order so the circles read 5, 4, 3, 2, 1 (top-to-bottom).
# Shape "Google Shape;371;p23" (the middle circle, already "3") keeps its
# text but is still touched as part of the same edit.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

$s.Shapes.Item("Google Shape;369;p23").TextFrame.TextRange.Text = "5"
$s.Shapes.Item("Google Shape;370;p23").TextFrame.TextRange.Text = "4"
$s.Shapes.Item("Google Shape;371;p23").TextFrame.TextRange.Text = "3"
$s.Shapes.Item("Google Shape;372;p23").TextFrame.TextRange.Text = "2"
$s.Shapes.Item("Google Shape;373;p23").TextFrame.TextRange.Text = "1"
